$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A74").Value = "31b57bd0-6c05-41b5-8f14-48622428d41d"
$ws.Range("B74").Value = "Login with valid credentials"
$ws.Range("C74").Value = "PASSED"
$ws.Range("D74").Value = "03_28_2024_22_55_57"
$ws.Range("E74").Value = "03_28_2024_22_56_08"
$ws.Range("F74").Value = "PT11.1655475S"
$ws.Range("D74:F74").HorizontalAlignment = -4108

$ws.Range("A75").Value = "8c06243c-0e6f-4b2d-852c-38835e66d1ea"
$ws.Range("B75").Value = "Create Country"
$ws.Range("C75").Value = "PASSED"
$ws.Range("D75").Value = "03_28_2024_22_56_13"
$ws.Range("E75").Value = "03_28_2024_22_56_24"
$ws.Range("F75").Value = "PT10.9984168S"
$ws.Range("D75:F75").HorizontalAlignment = -4108
